$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 4899
$ws.Range("C3").Value = 8031
$ws.Range("D3").Value = 3657
$ws.Range("E3").Value = 6553
$ws.Range("G3").Value = 12800
$ws.Range("B4").Value = 641.728512
$ws.Range("C4").Value = 1052.770304
$ws.Range("D4").Value = 479.199232
$ws.Range("E4").Value = 858.783744
$ws.Range("F4").Value = 1501.560832
$ws.Range("G4").Value = 1676.673024
$ws.Range("B5").Value = 203.04
$ws.Range("C5").Value = 241.44
$ws.Range("D5").Value = 1068.68
$ws.Range("E5").Value = 1154.35
$ws.Range("F5").Value = 1264.3
$ws.Range("G5").Value = 2046.16
$ws.Range("B6").Value = 383
$ws.Range("C6").Value = 404
$ws.Range("D6").Value = 1876
$ws.Range("E6").Value = 2057
$ws.Range("F6").Value = 2638
$ws.Range("G6").Value = 7701
$ws.Range("B7").Value = 396
$ws.Range("C7").Value = 433
$ws.Range("D7").Value = 2040
$ws.Range("E7").Value = 3163
$ws.Range("F7").Value = 3851
$ws.Range("G7").Value = 11076
$ws.Range("B12").Value = 4782
$ws.Range("C12").Value = 35500
$ws.Range("D12").Value = 98600
$ws.Range("E12").Value = 157000
$ws.Range("F12").Value = 60300
$ws.Range("G12").Value = 278000
$ws.Range("B13").Value = 19.6083712
$ws.Range("C13").Value = 145.752064
$ws.Range("D13").Value = 403.70176
$ws.Range("E13").Value = 641.728512
$ws.Range("F13").Value = 246.41536
$ws.Range("G13").Value = 1139.802112
$ws.Range("B14").Value = 207.89832
$ws.Range("C14").Value = 55.03463
$ws.Range("D14").Value = 39.57362000000001
$ws.Range("E14").Value = 39.81417
$ws.Range("F14").Value = 255.95412
$ws.Range("G14").Value = 84.10915
$ws.Range("C15").Value = 102.912
$ws.Range("D15").Value = 91.648
$ws.Range("E15").Value = 100.864
$ws.Range("F15").Value = 1548.288
$ws.Range("G15").Value = 432.128
$ws.Range("B16").Value = 1564.672
$ws.Range("C16").Value = 130.56
$ws.Range("D16").Value = 113.152
$ws.Range("E16").Value = 148.48
$ws.Range("F16").Value = 1859.584
$ws.Range("G16").Value = 626.688
$ws.Range("B21").Value = 6360
$ws.Range("C21").Value = 9846
$ws.Range("D21").Value = 14600
$ws.Range("E21").Value = 19500
$ws.Range("F21").Value = 17200
$ws.Range("G21").Value = 13000
$ws.Range("B22").Value = 833.61792
$ws.Range("C22").Value = 1290.797056
$ws.Range("D22").Value = 1917.845504
$ws.Range("E22").Value = 2550.136832
$ws.Range("F22").Value = 2251.292672
$ws.Range("G22").Value = 1835.008
$ws.Range("B23").Value = 77.03
$ws.Range("C23").Value = 85.53
$ws.Range("D23").Value = 102.28
$ws.Range("E23").Value = 128.51
$ws.Range("F23").Value = 180.28
$ws.Range("G23").Value = 1152.47
$ws.Range("B24").Value = 72
$ws.Range("C24").Value = 93
$ws.Range("D24").Value = 135
$ws.Range("E24").Value = 219
$ws.Range("F24").Value = 192
$ws.Range("G24").Value = 4015
$ws.Range("B25").Value = 217
$ws.Range("C25").Value = 1303
$ws.Range("D25").Value = 1029
$ws.Range("E25").Value = 1336
$ws.Range("F25").Value = 2147
$ws.Range("G25").Value = 11207
$ws.Range("B30").Value = 122000
$ws.Range("C30").Value = 165000
$ws.Range("D30").Value = 260000
$ws.Range("E30").Value = 392000
$ws.Range("F30").Value = 373000
$ws.Range("G30").Value = 251000
$ws.Range("B31").Value = 499.122176
$ws.Range("C31").Value = 677.380096
$ws.Range("D31").Value = 1063.256064
$ws.Range("E31").Value = 1607.467008
$ws.Range("F31").Value = 1527.775232
$ws.Range("G31").Value = 1028.653056
$ws.Range("B32").Value = 5.457020000000001
$ws.Range("F32").Value = 15.41
$ws.Range("G32").Value = 83.54000000000001
$ws.Range("B33").Value = 6.496
$ws.Range("D33").Value = 17
$ws.Range("E33").Value = 11
$ws.Range("G33").Value = 330
$ws.Range("B34").Value = 14.4
$ws.Range("C34").Value = 13
$ws.Range("E34").Value = 24
$ws.Range("F34").Value = 137
$ws.Range("G34").Value = 742
$ws.Range("B39").Value = 11800
$ws.Range("C39").Value = 21100
$ws.Range("D39").Value = 6390
$ws.Range("E39").Value = 7620
$ws.Range("F39").Value = 7772
$ws.Range("G39").Value = 8255
$ws.Range("B40").Value = 1542.455296
$ws.Range("C40").Value = 2767.192064
$ws.Range("D40").Value = 837.812224
$ws.Range("E40").Value = 999.292928
$ws.Range("F40").Value = 1019.215872
$ws.Range("G40").Value = 1082.130432
$ws.Range("B41").Value = 77.27
$ws.Range("C41").Value = 81.69
$ws.Range("D41").Value = 595.9299999999999
$ws.Range("E41").Value = 868.2
$ws.Range("F41").Value = 1647.09
$ws.Range("G41").Value = 3033.47
$ws.Range("B42").Value = 251
$ws.Range("C42").Value = 273
$ws.Range("D42").Value = 1926
$ws.Range("F42").Value = 5669
$ws.Range("G42").Value = 15664
$ws.Range("B43").Value = 562
$ws.Range("C43").Value = 453
$ws.Range("D43").Value = 2343
$ws.Range("E43").Value = 4948
$ws.Range("F43").Value = 15533
$ws.Range("G43").Value = 28443
$ws.Range("B48").Value = 144000
$ws.Range("C48").Value = 155000
$ws.Range("D48").Value = 215000
$ws.Range("E48").Value = 301000
$ws.Range("F48").Value = 267000
$ws.Range("G48").Value = 332000
$ws.Range("B49").Value = 588.251136
$ws.Range("C49").Value = 636.485632
$ws.Range("D49").Value = 879.755264
$ws.Range("E49").Value = 1231.028224
$ws.Range("F49").Value = 1092.616192
$ws.Range("G49").Value = 1360.003072
$ws.Range("B50").Value = 6.44061
$ws.Range("C50").Value = 10.20723
$ws.Range("D50").Value = 16.9717
$ws.Range("E50").Value = 24.61053
$ws.Range("F50").Value = 47.83172
$ws.Range("G50").Value = 88.0149
$ws.Range("B51").Value = 1.912
$ws.Range("C51").Value = 2.224
$ws.Range("D51").Value = 2.8
$ws.Range("E51").Value = 2.8
$ws.Range("F51").Value = 2.576
$ws.Range("B52").Value = 101.888
$ws.Range("C52").Value = 120.32
$ws.Range("D52").Value = 407.552
$ws.Range("E52").Value = 1220.608
$ws.Range("F52").Value = 1662.976
$ws.Range("G52").Value = 2039.808
$ws.Range("B57").Value = 6564
$ws.Range("C57").Value = 8427
$ws.Range("D57").Value = 9266
$ws.Range("G57").Value = 12900
$ws.Range("B58").Value = 860.880896
$ws.Range("C58").Value = 1104.150528
$ws.Range("D58").Value = 1214.251008
$ws.Range("E58").Value = 1642.070016
$ws.Range("F58").Value = 1732.247552
$ws.Range("G58").Value = 1686.110208
$ws.Range("B59").Value = 60.32
$ws.Range("C59").Value = 77.38
$ws.Range("D59").Value = 188
$ws.Range("E59").Value = 128.65
$ws.Range("F59").Value = 183.37
$ws.Range("G59").Value = 823.5700000000001
$ws.Range("B60").Value = 60
$ws.Range("D60").Value = 174
$ws.Range("E60").Value = 243
$ws.Range("F60").Value = 204
$ws.Range("G60").Value = 3359
$ws.Range("B61").Value = 281
$ws.Range("C61").Value = 105
$ws.Range("D61").Value = 245
$ws.Range("E61").Value = 420
$ws.Range("F61").Value = 2114
$ws.Range("G61").Value = 8291
$ws.Range("B66").Value = 137000
$ws.Range("C66").Value = 158000
$ws.Range("E66").Value = 334000
$ws.Range("F66").Value = 377000
$ws.Range("G66").Value = 239000
$ws.Range("B67").Value = 562.036736
$ws.Range("C67").Value = 644.87424
$ws.Range("E67").Value = 1367.343104
$ws.Range("F67").Value = 1544.552448
$ws.Range("G67").Value = 980.41856
$ws.Range("E68").Value = 12.26
$ws.Range("F68").Value = 10.03
$ws.Range("G68").Value = 87.39
$ws.Range("B69").Value = 5.536
$ws.Range("C69").Value = 10
$ws.Range("D69").Value = 12
$ws.Range("E69").Value = 24
$ws.Range("F69").Value = 10
$ws.Range("G69").Value = 293
$ws.Range("B70").Value = 23.68
$ws.Range("C70").Value = 11
$ws.Range("D70").Value = 15
$ws.Range("E70").Value = 32
$ws.Range("F70").Value = 18
$ws.Range("G70").Value = 848
